$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Pharmacology" to "Scanner"
$ws.Name = "Scanner"

# Convert the "Log Date" column (C2:C40) from DD/MM/YYYY to MM/DD/YYYY,
# keeping the values as literal text (not auto-converted to date serials).
# A leading apostrophe forces Excel to store it as text; resetting the
# style back to "Normal" removes the quote-prefix formatting flag that
# Excel would otherwise leave behind on the cell.
$ws.Range("C2:C40").Value = "'05/19/2025"
$ws.Range("C2:C40").Style = "Normal"
